# Add initial vendor data (Acrossvape, HussarVape, Vicious Ant, SvoëMesto, Unknown)
# and re-sort the "vendors" list alphabetically, matching the upstream commit
# "Common: Added some initial data for vendors, atomizers, mods".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vendors")
$ws.Activate()

# Append the 5 new vendor names right after the existing data (rows 2-60).
$ws.Range("A61").Value = "Acrossvape"
$ws.Range("A62").Value = "HussarVape"
$ws.Range("A63").Value = "Vicious Ant"
$ws.Range("A64").Value = "SvoëMesto"
$ws.Range("A65").Value = "Unknown"

# Re-sort the whole (now 64-row) vendor list alphabetically (ascending),
# matching the original sheet's sort order.
$dataRange = $ws.Range("A2:A65")
$keyRange = $ws.Range("A2:A65")
$dataRange.Sort($keyRange)

# Match the new selection recorded in the sheet view.
$ws.Range("A55").Select()

# Best-effort: restore the scrolled viewport recorded in the sheet/workbook views.
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
$win.Left = 1605
